# Update data in object to refer to filtered data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the filtered data values in row 10 (E10, F10); G10 is a shared
# formula (=SUM(E10:F10)) and will recalculate automatically.
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 3

# Move the active selection from A14 to D14.
$ws.Range("D14").Select()
